$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Feuil1")

# --- Swap the "dec"/"hex" header labels in B2/C2 --------------------------
# (column B now holds the hex input, column C the decimal result, matching
#  the new HEX2DEC direction of the conversion below)
$b2 = $ws.Range("B2").Value()
$c2 = $ws.Range("C2").Value()
$ws.Range("B2").Value = $c2
$ws.Range("C2").Value = $b2

# --- Refresh the centered alignment on the merged header cell B1:C1 -------
$ws.Range("B1:C1").HorizontalAlignment = -4108
$ws.Range("B1:C1").VerticalAlignment = -4108

# --- Flip the conversion column from DEC2HEX to HEX2DEC -------------------
# Row 3 keeps its own (non-shared) formula.
$ws.Range("C3").Formula = "=HEX2DEC(B3)"

# Rows 4-10 become one shared-formula block.
$ws.Range("C4:C10").Formula = "=HEX2DEC(B4)"

# Rows 12-14 become another shared-formula block (row 11 is a blank
# separator row between the two CAN-id groups).
$ws.Range("C12:C14").Formula = "=HEX2DEC(B12)"

# Row 16 (after the blank separator row 15) gets the same conversion.
$ws.Range("C16").Formula = "=HEX2DEC(B16)"

# The blank separator rows only carry the centered formatting, no value.
$ws.Range("C11").HorizontalAlignment = -4108
$ws.Range("C11").VerticalAlignment = -4108
$ws.Range("C15").HorizontalAlignment = -4108
$ws.Range("C15").VerticalAlignment = -4108

# --- Leave the selection on B3, matching the saved view --------------------
$ws.Range("B3").Select()
